$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.024.62'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').Value = '2.601.90'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '591.10'
$ws.Range('E5').Value = '  -2.19%  '
$ws.Range('D6').Value = '149.92'
$ws.Range('E6').Value = '  -3.26%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.75%  '
$ws.Range('D9').Value = '2.600.44'
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').Value = '0.128'
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('E13').Value = '  -3.39%  '
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('D15').Value = '3.070.62'
$ws.Range('E16').Value = '  -5.23%  '
$ws.Range('D17').Value = '66.900.52'
$ws.Range('E17').Value = '  -1.34%  '
$ws.Range('D18').Value = '2.600.72'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').Value = '364.66'
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('E20').Value = '  -2.09%  '
$ws.Range('D21').Value = '7.34'
$ws.Range('E21').Value = '  -4.67%  '
$ws.Range('D22').Value = '4.29'
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('D23').Value = '4.74'
$ws.Range('E23').Value = '  -4.89%  '
$ws.Range('D24').Value = '2.08'
$ws.Range('E24').Value = '  -1.39%  '
$ws.Range('E25').Value = '  +3.46%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').Value = '2.733.18'
$ws.Range('D29').Value = '584.36'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('D31').Value = '0.0₃0984'
$ws.Range('E31').Value = '  -7.09%  '
$ws.Range('E32').Value = '  -5.84%  '
$ws.Range('D33').Value = '7.65'
$ws.Range('E33').Value = '  -3.92%  '
$ws.Range('E34').Value = '  -3.60%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -5.38%  '
$ws.Range('E37').Value = '  -3.20%  '
$ws.Range('D38').Value = '155.78'
$ws.Range('E38').Value = '  -1.28%  '
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('E42').Value = '  -3.86%  '
$ws.Range('E43').Value = '  -4.28%  '
$ws.Range('D44').Value = '17.08'
$ws.Range('E44').Value = '  +3.91%  '
$ws.Range('D46').Value = '153.24'
$ws.Range('E46').Value = '  -2.81%  '
$ws.Range('E47').Value = '  -3.39%  '
$ws.Range('E48').Value = '  -1.69%  '
$ws.Range('E49').Value = '  -3.71%  '
$ws.Range('E50').Value = '  -1.66%  '
$ws.Range('D51').Value = '21.39'
$ws.Range('E51').Value = '  +1.23%  '
